$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new driver/placa records to the bottom of the table.
$newRows = @(
    @("CVP2519", "LUCAS RIBEIRO SANTOS", 619883),
    @("FKK4G59", "LUCAS RIBEIRO SANTOS", 619883),
    @("FVW5D39", "LUCAS RIBEIRO SANTOS", 619883),
    @("HSF-4E78", "LUCAS RIBEIRO SANTOS", 619883),
    @("FBB7E68", "RAFAEL WANDERLLEY NASCIMENTO DOS SANTOS", 619609)
)

$startRow = 66
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
